# Insert a new gene/disease row above the current row 18 (THBS1), shifting
# it (and everything below) down by one, then populate the new row and
# append one more row at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 18 (and below) down by inserting a new blank row at position 18.
$ws.Rows.Item(18).Insert()

# Fill the newly inserted row 18 with the SNCG entry.
$ws.Cells.Item(18, 1).Value = "SNCG"
$ws.Cells.Item(18, 2).Value = "Autosomal dominant Charcot-Marie-Tooth disease type 2L ORPHA:99945"

# Append the TIMP2 entry as the new last row (20).
$ws.Cells.Item(20, 1).Value = "TIMP2"
$ws.Cells.Item(20, 2).Value = "Charcot-Marie-Tooth disease type 1A ORPHA:101081, Charcot-Marie-Tooth disease type 1E ORPHA:90658"
